$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format before writing values so that
# numeric-looking strings (e.g. "1.00", "8.00") are preserved exactly as
# text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '68.210.93'
$ws.Range("E2").Value = '  -1.98%  '

# Row 3
$ws.Range("D3").Value = '2.444.63'
$ws.Range("E3").Value = '  -1.99%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '553.35'
$ws.Range("E5").Value = '  -2.76%  '

# Row 6
$ws.Range("D6").Value = '159.64'
$ws.Range("E6").Value = '  -3.28%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '0.498'
$ws.Range("E8").Value = '  -2.69%  '

# Row 9
$ws.Range("D9").Value = '2.444.13'
$ws.Range("E9").Value = '  -1.93%  '

# Row 10
$ws.Range("E10").Value = '  -7.78%  '

# Row 11
$ws.Range("E11").Value = '  -1.37%  '

# Row 12
$ws.Range("E12").Value = '  -6.29%  '

# Row 13
$ws.Range("E13").Value = '  -3.60%  '

# Row 14
$ws.Range("D14").Value = '2.895.45'
$ws.Range("E14").Value = '  -1.68%  '

# Row 15
$ws.Range("D15").Value = '68.042.14'
$ws.Range("E15").Value = '  -2.01%  '

# Row 16
$ws.Range("E16").Value = '  -5.55%  '

# Row 17
$ws.Range("D17").Value = '23.05'
$ws.Range("E17").Value = '  -5.17%  '

# Row 18
$ws.Range("D18").Value = '2.440.76'
$ws.Range("E18").Value = '  -1.34%  '

# Row 19
$ws.Range("D19").Value = '10.66'
$ws.Range("E19").Value = '  -4.61%  '

# Row 20
$ws.Range("D20").Value = '338.28'
$ws.Range("E20").Value = '  -2.31%  '

# Row 21
$ws.Range("E21").Value = '  -5.72%  '

# Row 22
$ws.Range("D22").Value = '3.74'
$ws.Range("E22").Value = '  -3.55%  '

# Row 23
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.04%  '

# Row 24
$ws.Range("E24").Value = '  -3.78%  '

# Row 25
$ws.Range("D25").Value = '66.04'
$ws.Range("E25").Value = '  -5.21%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -7.25%  '

# Row 27
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.571.37'
$ws.Range("E27").Value = '  -1.70%  '

# Row 28
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.01%  '

# Row 29
$ws.Range("D29").Value = '8.00'
$ws.Range("E29").Value = '  -7.80%  '

# Row 30
$ws.Range("E30").Value = '  -8.47%  '

# Row 31
$ws.Range("D31").Value = '7.06'
$ws.Range("E31").Value = '  -7.81%  '

# Row 32
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.03%  '

# Row 33
$ws.Range("D33").Value = '426.29'
$ws.Range("E33").Value = '  -3.01%  '

# Row 34
$ws.Range("E34").Value = '  -6.19%  '

# Row 35
$ws.Range("E35").Value = '  -6.09%  '

# Row 36
$ws.Range("D36").Value = '154.96'
$ws.Range("E36").Value = '  -0.06%  '

# Row 37
$ws.Range("D37").Value = '19.01'
$ws.Range("E37").Value = '  -0.28%  '

# Row 38
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("E39").Value = '  -4.05%  '

# Row 40
$ws.Range("E40").Value = '  -2.86%  '

# Row 41
$ws.Range("E41").Value = '  -4.96%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '4.33'
$ws.Range("E42").Value = '  -5.64%  '

# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '37.28'
$ws.Range("E43").Value = '  -1.62%  '

# Row 44
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '1.09'
$ws.Range("E44").Value = '  +1.25%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.45'
$ws.Range("E45").Value = '  -8.63%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '1.99'
$ws.Range("E46").Value = '  -8.11%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '130.90'
$ws.Range("E47").Value = '  -5.55%  '

# Row 48
$ws.Range("D48").Value = '3.30'
$ws.Range("E48").Value = '  -4.08%  '

# Row 49
$ws.Range("E49").Value = '  -1.97%  '

# Row 50
$ws.Range("E50").Value = '  -7.38%  '

# Row 51
$ws.Range("D51").Value = '0.556'
$ws.Range("E51").Value = '  -2.96%  '

# Restore the original (unset) number format so the cells serialize
# back out exactly like the rest of the sheet (no stray style refs).
$ws.Range("D2:D51").ClearFormats()
